$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}

Replace-Text "2024-10-23 Wednesday" "2024-10-24 Thursday"

Replace-Text "83×65=" "74×67="
Replace-Text "77×95=" "21×26="
Replace-Text "63×85=" "93×33="
Replace-Text "86×76=" "66×88="
Replace-Text "27×33=" "18×17="

Replace-Text "92×19=" "80×99="
Replace-Text "39×64=" "62×85="
Replace-Text "62×52=" "85×81="
Replace-Text "93×28=" "23×15="
Replace-Text "74×54=" "49×11="

Replace-Text "33×11=" "96×79="
Replace-Text "31×24=" "52×85="
Replace-Text "27×15=" "95×67="
Replace-Text "23×49=" "16×25="
Replace-Text "15×87=" "43×56="

Replace-Text "35×73=" "88×31="
Replace-Text "21×50=" "42×97="
Replace-Text "65×64=" "94×49="
Replace-Text "60×38=" "38×21="
Replace-Text "56×70=" "13×80="

Replace-Text "66×39=" "80×51="
Replace-Text "12×50=" "40×16="
Replace-Text "67×99=" "45×28="
Replace-Text "70×60=" "58×96="
Replace-Text "90×42=" "27×44="
